$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# Update URL value (row 2)
$ws1.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/vaccine-gender"

# Update Date value (row 8)
$ws1.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (row 11).
# Shift existing rows 11-14 down to 12-15 by copying content bottom-up so we don't
# clobber values before they're copied. This also keeps the existing per-cell style
# (s="2") intact on all the rows that already existed.
for ($r = 14; $r -ge 11; $r--) {
    $srcA = $ws1.Cells.Item($r, 1)
    $srcB = $ws1.Cells.Item($r, 2)
    $dstA = $ws1.Cells.Item($r + 1, 1)
    $dstB = $ws1.Cells.Item($r + 1, 2)
    $dstA.Value = $srcA.Text
    $dstB.Value = $srcB.Text
}

# The sheet grew by one row (A1:B14 -> A1:B15); row 15 is brand new, so it has no
# style yet. Copy the row-14 formatting down onto it before putting its real values in.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Now fill in the new "Jurisdiction" row (row 11) with blank value
$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""

# --- Sheet 2: "Include from Vaccine Gender" -> "Include #0" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# Update System URI value (row 7)
$ws2.Cells.Item(7, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/VaccineGender"
